$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add row 119: 11/26/2020 (serial 44161), Tests=0, Positive=0
$ws.Range("A119").Value = 44161
$ws.Range("A119").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"
$ws.Range("B119").Value = 0
$ws.Range("C119").Value = 0

# Add row 120: 11/27/2020 (serial 44162), Tests=0, Positive=0
$ws.Range("A120").Value = 44162
$ws.Range("A120").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"
$ws.Range("B120").Value = 0
$ws.Range("C120").Value = 0

# Update selection to F117 as per diff
$ws.Range("F117").Select()
